# RCC new Script implementation
# Adds a new "RCC113" watch-list test case as row 15 on the "Test Cases" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Values --------------------------------------------------------------
# Write Description (C) before Jira id (B) so the shared-string table picks
# up the same insertion order as the rest of the sheet (TCID, Description,
# Jira id) - e.g. rows 5/12 follow that same pattern.
$ws.Range("A15").Value = "RCC113"
$ws.Range("C15").Value = "Verify that user is able to add an article to the group from watch list details page. ||Verify that user is able to add a post to the group from watch list details page.||Verify that user is able to add a patent to the group from watch list details page."
$ws.Range("B15").Value = "OPQA-3456||OPQA-3460||OPQA-3464"
$ws.Range("D15").Value = "Y"
$ws.Range("E15").Value = ""

# --- Formatting ------------------------------------------------------------
# A15 / D15 / E15: plain bordered cells, same look as the rest of column A/D/E.
$ws.Range("A2").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = "RCC113"

$ws.Range("A2").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = "Y"

$ws.Range("A2").Copy()
$ws.Range("E15").PasteSpecial(-4122)

# B15: bordered + wrap text, matching the other Jira-id cells that hold
# pipe-separated lists (e.g. B12/B13).
$ws.Range("B12").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = "OPQA-3456||OPQA-3460||OPQA-3464"

# C15: new style - bordered, left/top aligned, wrap text, with a light
# (theme "Background 1") fill so the long description stands out.
$ws.Range("C15").WrapText = $true
$ws.Range("C15").HorizontalAlignment = -4131
$ws.Range("C15").VerticalAlignment = -4160
$ws.Range("C15").Interior.Pattern = 1
$ws.Range("C15").Interior.ThemeColor = 2
$ws.Range("C15").Interior.PatternColorIndex = -4105

# Row 15 is a tall, wrapped row.
$ws.Rows.Item(15).RowHeight = 45

# --- Selection --------------------------------------------------------------
$ws.Range("C18").Select()
